# Corrected image of lab compendium
#
# Removes four stray shapes that were mistakenly left on the slide:
#   - "Rectangle 63"              (small anchor rectangle)
#   - "Connecteur : en angle 64"  (bent connector arrow anchored to the rectangle)
#   - "ZoneTexte 65"              ("0..*" text label)
#   - "ZoneTexte 66"              ("action.(...).definitionCanonical (reflexes)" text label)
#
# These shapes are identified by their stable PowerPoint shape Id
# (rather than by Name) to avoid any ambiguity with special/non-ASCII
# characters that PowerPoint inserts into some autogenerated shape names.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$idsToRemove = @(64, 65, 66, 67)

foreach ($targetId in $idsToRemove) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Id -eq $targetId) {
            $sh.Delete()
            break
        }
    }
}
